$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells for the two new columns ---
$ws.Range("H1").Value = "insert string"
$ws.Range("I1").Value = "values string"

# --- Column widths for the new H and I columns ---
$ws.Columns.Item(8).ColumnWidth = 206.66666666666666
$ws.Columns.Item(9).ColumnWidth = 53.498697916666664

# --- Fill column I (rows 2-51) with the VALUES-array formula ---
for ($r = 2; $r -le 51; $r++) {
    $formula = '=CONCATENATE("[",A' + $r + ',",""",B' + $r + ',""",""",C' + $r + ',""",""",D' + $r + ',""",""",E' + $r + ',""",""",F' + $r + ',""",""",G' + $r + ',"""],")'
    $ws.Range("I$r").Formula = $formula
}

# --- Update the active selection/view to match the new focus on column H/I ---
$ws.Range("H10").Select()
